$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B (rows 2-48, corresponding to A=0..46)
$bValues = @(
    0.8497306202089479,
    2.935541320861634,
    3.307866141770259,
    5.701211927436302,
    6.957028089973347,
    7.351241215549599,
    7.920421814630269,
    12.4980082144637,
    13.61541676083227,
    15.94428146348539,
    19.13922303899406,
    21.28099878412854,
    25.50485956133809,
    26.51076955105791,
    34.60021111813396,
    36.43006116509981,
    39.18376976365121,
    39.71185182635229,
    39.95885669990316,
    40.1864458979371,
    41.49955496993493,
    44.39813841966532,
    47.30263598979391,
    47.78698397867332,
    49.4984593930864,
    49.89796625862768,
    50.92223511696007,
    54.95704749794656,
    58.63072393140954,
    60.31058917854543,
    60.54277840515179,
    63.94802032991244,
    69.64055892116387,
    70.03112457238136,
    72.56042996476505,
    76.34877935700437,
    76.88400044418229,
    78.50492656357167,
    78.63324987084519,
    79.21785093293826,
    80.9219025823228,
    88.35183388092935,
    89.14842979410734,
    92.28113899908476,
    96.84986321480703,
    97.62277343948121,
    99.33075746032056
)

# New values for column C (rows 2-48, corresponding to A=0..46)
$cValues = @(
    1.66993850619004,
    3.319242904023441,
    4.769167891355981,
    6.286812531847147,
    7.371866663776511,
    8.962268866839073,
    10.66089065248301,
    12.4018605283593,
    13.62113227475125,
    15.13855134020944,
    16.29579942229515,
    18.25744292921607,
    19.60040237699198,
    21.13696139515462,
    22.47227450741326,
    23.77255391857918,
    25.06589665113383,
    26.19816686124863,
    27.76145496906675,
    29.55305704522209,
    30.79153583739203,
    32.10976988238109,
    34.05478676198223,
    35.4625986237866,
    37.08562966033011,
    38.69508024229678,
    40.37874156138224,
    41.95140332576854,
    43.24415665142383,
    44.8572829353081,
    46.53255045599396,
    49.57652659226678,
    50.82350461759883,
    52.3896429630467,
    54.40845637139464,
    56.13798793569574,
    57.31941497915404,
    59.0435271410422,
    60.29019684164526,
    61.81449091057284,
    63.40730177004358,
    64.61219045535587,
    66.42956696734723,
    67.65419933846869,
    69.04874559135781,
    70.5733000991241,
    72.08403712205855
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

# Remove the now-obsolete last row (previously row 49, A=47)
$ws.Rows.Item(49).Delete()
